$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 34.65554307931001
$ws.Range("G2").Value = 33.24135261438737
$ws.Range("H2").Value = 36.04062058318276
$ws.Range("I2").Value = 3.469365924523141
$ws.Range("J2").Value = 3.208344111345506
$ws.Range("K2").Value = 3.733288334860029
$ws.Range("L2").Value = 0.254592384888277
$ws.Range("M2").Value = 0.235945042664084
$ws.Range("N2").Value = 0.2732372777061016

# Row 3
$ws.Range("F3").Value = 0.005082095168829363
$ws.Range("G3").Value = 0.003402116081251368
$ws.Range("H3").Value = 0.007055646502701864
$ws.Range("I3").Value = 0.004653804487702032
$ws.Range("J3").Value = 0.003104472924203355
$ws.Range("K3").Value = 0.006481977647927694
$ws.Range("L3").Value = 0.005043944123421878
$ws.Range("M3").Value = 0.003371958872956345
$ws.Range("N3").Value = 0.007010814451036995

# Row 4
$ws.Range("F4").Value = 34.66062517447885
$ws.Range("G4").Value = 33.24475473046861
$ws.Range("H4").Value = 36.04767622968546
$ws.Range("I4").Value = 3.474019729010843
$ws.Range("J4").Value = 3.21144858426971
$ws.Range("K4").Value = 3.739770312507956
$ws.Range("L4").Value = 0.2596363290116989
$ws.Range("M4").Value = 0.2393170015370403
$ws.Range("N4").Value = 0.2802480921571386
